# Update crypto price/volume figures on the active worksheet (data refresh
# performed by the "Updated symbol list" GitHub Actions job).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell reference -> new text value. These are stored as plain text
# in the sheet (prices/percentages with fixed decimal places), so each
# value is written while the cell is temporarily forced to Text format,
# then the format override is cleared again to avoid changing the cell's
# style from its original (unstyled) state.
$changes = [ordered]@{
    "D2"  = "261.05"
    "E2"  = "1.78%"
    "D3"  = "27.39"
    "E3"  = "1.78%"
    "D4"  = "4.703"
    "E4"  = "-2.07%"
    "D5"  = "0.06090"
    "E5"  = "2.04%"
    "D6"  = "6.674"
    "E6"  = "0.55%"
    "D7"  = "0.8460"
    "D8"  = "0.9290"
    "E8"  = "0.34%"
    "E9"  = "1.93%"
    "D10" = "0.04832"
    "E10" = "14.84%"
    "D11" = "0.07102"
    "D12" = "0.03086"
    "E12" = "1.02%"
    "D13" = "0.09061"
    "D14" = "0.001542"
    "E14" = "0.82%"
    "D15" = "0.0006104"
    "D16" = "0.006032"
    "E16" = "-0.83%"
    "E17" = "-0.60%"
    "D18" = "3.149"
    "E18" = "-0.31%"
    "E20" = "2.25%"
    "E21" = "0.21%"
    "D22" = "4.109"
    "E22" = "5.26%"
    "D23" = "0.04255"
    "E23" = "0.19%"
    "E24" = "0.16%"
    "D25" = "0.003793"
    "E25" = "4.84%"
    "D26" = "0.0001200"
    "E26" = "0.05%"
    "E27" = "3.46%"
    "D40" = "0.03871"
    "E40" = "2.34%"
    "E41" = "1.32%"
    "E42" = "-34.21%"
    "D43" = "0.01623"
    "E43" = "14.82%"
    "E44" = "-9.49%"
    "D45" = "0.00005142"
    "E45" = "-3.45%"
    "E46" = "0.09%"
    "E47" = "-39.23%"
    "E48" = "23.82%"
    "E49" = "0.09%"
    "E50" = "0.09%"
}

foreach ($ref in $changes.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $changes[$ref]
    $cell.ClearFormats()
}
